# Workbook has sheets (in tab order): Bus, Load, Gen slack, Lines, Trans
$wb = $excel.ActiveWorkbook

# --- Content change -------------------------------------------------
# Trans!B2 std_type label: "0.25 MVA 10/0.4 kV" -> "0.25 MVA 20/0.4 kV"
$transSheet = $wb.Worksheets.Item("Trans")
$transSheet.Activate()
$transSheet.Range("B2").Value = "0.25 MVA 20/0.4 kV"
$transSheet.Range("F10").Select()

# --- Selection / active-tab changes ---------------------------------
# Bus sheet: move the cell selection to B2
$busSheet = $wb.Worksheets.Item("Bus")
$busSheet.Activate()
$busSheet.Range("B2").Select()

# Load sheet: becomes the active sheet/tab (selection stays at F14)
$loadSheet = $wb.Worksheets.Item("Load")
$loadSheet.Activate()
